# Add newly annotated rows to the "Non OCRED - OCR Problems" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Non OCRED - OCR Problems")

# Row 24
$ws.Range("A24").Value = 3
$ws.Range("A24").HorizontalAlignment = -4152
$ws.Range("B24").Value = "https://demo.humlab.umu.se/courier/081400engo.pdf"
$ws.Range("C24").Value = "Contrast"

# Row 25
$ws.Range("A25").Value = 5
$ws.Range("A25").HorizontalAlignment = -4152
$ws.Range("B25").Value = "https://demo.humlab.umu.se/courier/081370engo.pdf"
$ws.Range("C25").Value = "Contrast"

# Row 26
$ws.Range("A26").Value = "2,7,9,12"
$ws.Range("A26").HorizontalAlignment = -4152
$ws.Range("B26").Value = "https://demo.humlab.umu.se/courier/074977engo.pdf"
$ws.Range("C26").Value = "Contrast"

# Row 27
$ws.Range("A27").Value = "3,5,9,10"
$ws.Range("A27").HorizontalAlignment = -4152
$ws.Range("B27").Value = "https://demo.humlab.umu.se/courier/074946engo.pdf"
$ws.Range("C27").Value = "Contrast, full article. Whole issue has some trouble when italic font, and with contrast"

# Row 28
$ws.Range("A28").Value = 8
$ws.Range("A28").HorizontalAlignment = -4152
$ws.Range("B28").Value = "https://demo.humlab.umu.se/courier/074184engo.pdf"
$ws.Range("C28").Value = "Article"

# Move the active selection, matching the author's last cursor position.
$ws.Activate()
$ws.Range("P18").Select()
